$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text cell whose content happens to look like a
# boolean keyword ("true"/"false") without Excel auto-converting it to a
# real Boolean cell type and without leaving a stray quote-prefixed style
# behind. Round-tripping through a TEXT() formula then pasting values-only
# keeps it a genuine shared-string text cell.
function Set-LiteralText($range, [string]$text) {
    $escaped = $text.Replace('"', '""""')
    $range.Formula = '=TEXT(1,"""' + $escaped + '""")'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# --- Header row (row 1) --- (establishes the shared-string order for the
# column headers first, matching how the source data was serialised)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "sib.name"
$ws.Range("D1").Value = "sib.ph"
$ws.Range("E1").Value = "sib.addr"
$ws.Range("F1").Value = "frnds.b"
$ws.Range("G1").Value = "frnds.best"

# --- Row 2 ---
$ws.Range("B2").Value = "aditi"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 98765
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = 3

# --- Row 3 (re-uses "null") ---
$ws.Range("B3").Clear()
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 12345
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = "null"

# --- Row 4 (introduces "fjaslkff kjas" before the boolean-looking text) ---
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "fjaslkff kjas"
$ws.Range("F4").Value = 3

# --- Now the boolean-looking literals, in row order ---
Set-LiteralText $ws.Range("G2") "true"
Set-LiteralText $ws.Range("G4") "false"

# Copy the bordered/bold "index column" formatting from A4 down onto the
# newly added A5:A7 cells (A2:A4 already carry it from the source sheet).
$ws.Range("A4").Copy()
$ws.Range("A5:A7").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 5 ---
$ws.Range("A5").Value = 3
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = "null"

# --- Row 6 ---
$ws.Range("A6").Value = 4
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = "null"

# --- Row 7 ---
$ws.Range("A7").Value = 5
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = "null"

# --- Remove the old H:J columns entirely (age/secretIdentity/powers data no longer used) ---
$ws.Range("H1:J4").Clear()
